$d = $word.ActiveDocument

$pairs = @(
    @("610÷6=101, 4", "619÷7=88, 3"),
    @("232÷8=29, 0", "291÷6=48, 3"),
    @("743÷6=123, 5", "957÷6=159, 3"),
    @("989÷7=141, 2", "171÷7=24, 3"),
    @("895÷7=127, 6", "405÷9=45, 0"),
    @("625÷5=125, 0", "776÷8=97, 0"),
    @("663÷2=331, 1", "875÷4=218, 3"),
    @("640÷4=160, 0", "161÷6=26, 5"),
    @("135÷8=16, 7", "541÷3=180, 1"),
    @("910÷6=151, 4", "724÷2=362, 0"),
    @("140÷9=15, 5", "489÷2=244, 1"),
    @("242÷6=40, 2", "584÷5=116, 4"),
    @("653÷6=108, 5", "983÷6=163, 5"),
    @("951÷8=118, 7", "452÷2=226, 0"),
    @("826÷4=206, 2", "618÷9=68, 6"),
    @("318÷6=53, 0", "142÷6=23, 4"),
    @("230÷3=76, 2", "708÷5=141, 3"),
    @("102÷9=11, 3", "808÷2=404, 0"),
    @("695÷9=77, 2", "379÷4=94, 3"),
    @("793÷8=99, 1", "582÷5=116, 2"),
    @("863÷4=215, 3", "456÷4=114, 0"),
    @("732÷4=183, 0", "307÷7=43, 6"),
    @("505÷8=63, 1", "854÷5=170, 4"),
    @("567÷5=113, 2", "900÷4=225, 0"),
    @("495÷7=70, 5", "183÷5=36, 3")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    } else {
        Write-Host "Replaced: $old -> $new"
    }
}
